$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[ S. M.%Ahmed%null%1,                           R. U.% Shah%null%1,                           M.% Bale%null%1,                           J. B.% Peacock%null%1,                           B.% Berger%null%1,                           A.% Brown%null%1,                           S.% Mann%null%1,                           W.% West%null%1,                           V.% Martin%null%1,                           V.% Fernandez%null%1,                           S.% Grineski%null%1,                           B. J.% Brintz%null%1,                           M. H.% Samore%null%1,                           M. J.% Ferrari%null%1,                           D. T.% Leung%null%1,                           L. T. % Keegan%null%1,         S. M.%Ahmed%null%1,         R. U.% Shah%null%1,         M.% Bale%null%1,         J. B.% Peacock%null%1,         B.% Berger%null%1,         A.% Brown%null%1,         S.% Mann%null%1,         W.% West%null%1,         V.% Martin%null%1,         V.% Fernandez%null%1,         S.% Grineski%null%1,         B. J.% Brintz%null%1,         M. H.% Samore%null%1,         M. J.% Ferrari%null%1,         D. T.% Leung%null%1,         L. T. % Keegan%null%1]'
$ws.Range('E3').Value = '[Michael G%Argenziano%NULL%5,                          Samuel L%Bruce%NULL%1,                          Cody L%Slater%NULL%1,                          Jonathan R%Tiao%NULL%1,                          Matthew R%Baldwin%NULL%1,                          R Graham%Barr%NULL%1,                          Bernard P%Chang%NULL%1,                          Katherine H%Chau%NULL%1,                          Justin J%Choi%NULL%1,                          Nicholas%Gavin%NULL%1,                          Parag%Goyal%NULL%1,                          Angela M%Mills%NULL%1,                          Ashmi A%Patel%NULL%1,                          Marie-Laure S%Romney%NULL%1,                          Monika M%Safford%NULL%1,                          Neil W%Schluger%NULL%1,                          Soumitra%Sengupta%NULL%1,                          Magdalena E%Sobieszczyk%NULL%1,                          Jason E%Zucker%NULL%1,                          Paul A%Asadourian%NULL%1,                          Fletcher M%Bell%NULL%1,                          Rebekah%Boyd%NULL%1,                          Matthew F%Cohen%NULL%1,                          MacAlistair I%Colquhoun%NULL%1,                          Lucy A%Colville%NULL%1,                          Joseph H%de Jonge%NULL%1,                          Lyle B%Dershowitz%NULL%1,                          Shirin A%Dey%NULL%1,                          Katherine A%Eiseman%NULL%1,                          Zachary P%Girvin%NULL%1,                          Daniella T%Goni%NULL%1,                          Amro A%Harb%NULL%1,                          Nicholas%Herzik%NULL%1,                          Sarah%Householder%NULL%1,                          Lara E%Karaaslan%NULL%1,                          Heather%Lee%NULL%1,                          Evan%Lieberman%NULL%1,                          Andrew%Ling%NULL%1,                          Ree%Lu%NULL%1,                          Arthur Y%Shou%NULL%1,                          Alexander C%Sisti%NULL%1,                          Zachary E%Snow%NULL%1,                          Colin P%Sperring%NULL%1,                          Yuqing%Xiong%NULL%1,                          Henry W%Zhou%NULL%1,                          Karthik%Natarajan%NULL%1,                          George%Hripcsak%NULL%1,                          Ruijun%Chen%NULL%1]'
$ws.Range('E4').Value = '[César%Caraballo%NULL%2,                          Megan%McCullough%NULL%4,                          Megan%McCullough%NULL%0,                          Michael A.%Fuery%NULL%2,                          Fouad%Chouairi%NULL%4,                          Fouad%Chouairi%NULL%0,                          Craig%Keating%NULL%2,                          Neal G.%Ravindra%NULL%2,                          P. Elliott%Miller%NULL%2,                          Maricar%Malinis%NULL%2,                          Nitu%Kashyap%NULL%2,                          Allen%Hsiao%NULL%2,                          F. Perry%Wilson%NULL%2,                          Jeptha P.%Curtis%NULL%2,                          Matthew%Grant%NULL%2,                          Eric J.%Velazquez%NULL%2,                          Nihar R.%Desai%NULL%2,                          Tariq%Ahmad%NULL%2,                          Chiara%Lazzeri%NULL%18,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0]'
$ws.Range('E5').Value = '[Jeremy A. W.%Gold%NULL%2,                          Karen K.%Wong%NULL%1,                          Christine M.%Szablewski%NULL%1,                          Priti R.%Patel%NULL%1,                          John%Rossow%NULL%1,                          Juliana%da Silva%NULL%1,                          Pavithra%Natarajan%NULL%1,                          Sapna Bamrah%Morris%NULL%1,                          Robyn Neblett%Fanfair%NULL%1,                          Jessica%Rogers-Brown%NULL%1,                          Beau B.%Bruce%NULL%1,                          Sean D.%Browning%NULL%1,                          Alfonso C.%Hernandez-Romieu%NULL%2,                          Nathan W.%Furukawa%NULL%1,                          Mohleen%Kang%NULL%1,                          Mary E.%Evans%NULL%1,                          Nadine%Oosmanally%NULL%1,                          Melissa%Tobin-D’Angelo%NULL%1,                          Cherie%Drenzek%NULL%1,                          David J.%Murphy%NULL%0,                          Julie%Hollberg%NULL%1,                          James M.%Blum%NULL%0,                          Robert%Jansen%NULL%1,                          David W.%Wright%NULL%1,                          William M.%Sewell%NULL%1,                          Jack D.%Owens%NULL%1,                          Benjamin%Lefkove%NULL%1,                          Frank W.%Brown%NULL%1,                          Deron C.%Burton%NULL%1,                          Timothy M.%Uyeki%NULL%3,                          Stephanie R.%Bialek%NULL%1,                          Brendan R.%Jackson%NULL%1]'
$ws.Range('E6').Value = '[Tian%Gu%NULL%1,                          Jasmine A.%Mack%NULL%1,                          Maxwell%Salvatore%NULL%1,                          Swaraaj Prabhu%Sankar%NULL%1,                          Thomas S.%Valley%NULL%1,                          Karandeep%Singh%NULL%1,                          Brahmajee K.%Nallamothu%NULL%1,                          Sachin%Kheterpal%NULL%1,                          Lynda%Lisabeth%NULL%1,                          Lars G.%Fritsche%NULL%1,                          Bhramar%Mukherjee%NULL%1,                           T.%Gu%null%1,                           J. A.% Mack%null%1,                           M.% Salvatore%null%1,                           S. P.% Sankar%null%1,                           T. S.% Valley%null%1,                           K.% Singh%null%1,                           B. K.% Nallamothu%null%1,                           S.% Kheterpal%null%1,                           L.% Lisabeth%null%1,                           L. G.% Fritsche%null%1,                           B. G. % Mukherjee%null%1,                          T.%Gu%null%1,                          J. A.% Mack%null%1,                          M.% Salvatore%null%1,                          S. P.% Sankar%null%1,                          T. S.% Valley%null%1,                          K.% Singh%null%1,                          B. K.% Nallamothu%null%1,                          S.% Kheterpal%null%1,                          L.% Lisabeth%null%1,                          L. G.% Fritsche%null%1,                          B. G. % Mukherjee%null%1,         T.%Gu%null%1,         J. A.% Mack%null%1,         M.% Salvatore%null%1,         S. P.% Sankar%null%1,         T. S.% Valley%null%1,         K.% Singh%null%1,         B. K.% Nallamothu%null%1,         S.% Kheterpal%null%1,         L.% Lisabeth%null%1,         L. G.% Fritsche%null%1,         B. G. % Mukherjee%null%1]'
$ws.Range('E7').Value = '[ T.%Jun%null%1,                           S.% Nirenberg%null%1,                           P.% Kovatch%null%1,                           K.-l. % Huang%null%1,         T.%Jun%null%1,         S.% Nirenberg%null%1,         P.% Kovatch%null%2,         K.-l. % Huang%null%1]'
$ws.Range('E8').Value = '[Lindsay%Kim%LKim@cdc.gov%1,                          Shikha%Garg%NULL%1,                          Alissa%O’Halloran%NULL%1,                          Michael%Whitaker%NULL%1,                          Huong%Pham%NULL%1,                          Evan J%Anderson%NULL%1,                          Isaac%Armistead%NULL%1,                          Nancy M%Bennett%NULL%1,                          Laurie%Billing%NULL%1,                          Kathryn%Como-Sabetti%NULL%1,                          Mary%Hill%NULL%1,                          Sue%Kim%NULL%1,                          Maya L%Monroe%NULL%1,                          Alison%Muse%NULL%1,                          Arthur L%Reingold%NULL%1,                          William%Schaffner%NULL%1,                          Melissa%Sutton%NULL%1,                          H Keipp%Talbot%NULL%1,                          Salina M%Torres%NULL%1,                          Kimberly%Yousey-Hindes%NULL%1,                          Rachel%Holstein%NULL%1,                          Charisse%Cummings%NULL%1,                          Lynette%Brammer%NULL%1,                          Aron J%Hall%NULL%1,                          Alicia M%Fry%NULL%1,                          Gayle E%Langley%NULL%1]'
$ws.Range('E9').Value = '[ T. J.%Levy%null%1,                           S.% Richardson%null%1,                           K.% Coppa%null%1,                           D. P.% Barnaby%null%1,                           T.% McGinn%null%1,                           L. B.% Becker%null%1,                           K. W.% Davidson%null%1,                           S. L.% Cohen%null%1,                           J. S.% Hirsch%null%1,                           T.% Zanos%null%1,                             % Northwell & Maimonides COVID-19 Research Consortium%null%1,         T. J.%Levy%null%1,         S.% Richardson%null%1,         K.% Coppa%null%1,         D. P.% Barnaby%null%1,         T.% McGinn%null%1,         L. B.% Becker%null%1,         K. W.% Davidson%null%1,         S. L.% Cohen%null%1,         J. S.% Hirsch%null%1,         T.% Zanos%null%1,           % Northwell & Maimonides COVID-19 Research Consortium%null%1]'
$ws.Range('E10').Value = '[ C.-H.%Lo%null%1,                           L. H.% Nguyen%null%1,                           D. A.% Drew%null%1,                           M. S.% Graham%null%1,                           E. T.% Warner%null%1,                           A. D.% Joshi%null%1,                           C. M.% Astley%null%1,                           C.-G.% Guo%null%1,                           W.% Ma%null%1,                           R. S.% Mehta%null%1,                           S.% Kwon%null%1,                           M.% Song%null%1,                           R.% Davies%null%1,                           J.% Capdevila%null%1,                           K. A.% Lee%null%1,                           M. N.% Lochlainn%null%1,                           T.% Varsavsky%null%1,                           C. H.% Sudre%null%1,                           J.% Wolf%null%1,                           Y. C.% Cozier%null%1,                           L.% Rosenberg%null%1,                           L. R.% Wilkens%null%1,                           C. A.% Haiman%null%1,                           L. L.% Marchand%null%1,                           J. R.% Palmer%null%1,                           T. D.% Spector%null%1,                           S.% Ourselin%null%1,                           C. J.% Steves%null%1,                           A. T.% Chan%null%1,                             % COPE Consortium%null%1,         C.-H.%Lo%null%1,         L. H.% Nguyen%null%1,         D. A.% Drew%null%1,         M. S.% Graham%null%1,         E. T.% Warner%null%1,         A. D.% Joshi%null%1,         C. M.% Astley%null%1,         C.-G.% Guo%null%1,         W.% Ma%null%1,         R. S.% Mehta%null%1,         S.% Kwon%null%1,         M.% Song%null%1,         R.% Davies%null%1,         J.% Capdevila%null%1,         K. A.% Lee%null%1,         M. N.% Lochlainn%null%1,         T.% Varsavsky%null%1,         C. H.% Sudre%null%1,         J.% Wolf%null%1,         Y. C.% Cozier%null%1,         L.% Rosenberg%null%1,         L. R.% Wilkens%null%1,         C. A.% Haiman%null%1,         L. L.% Marchand%null%1,         J. R.% Palmer%null%1,         T. D.% Spector%null%1,         S.% Ourselin%null%1,         C. J.% Steves%null%1,         A. T.% Chan%null%1,           % COPE Consortium%null%1]'
$ws.Range('E11').Value = '[Simon%de Lusignan%simon.delusignan@phc.ox.ac.uk%1,                          Jienchi%Dorward%NULL%1,                          Ana%Correa%NULL%1,                          Nicholas%Jones%NULL%1,                          Oluwafunmi%Akinyemi%NULL%1,                          Gayatri%Amirthalingam%NULL%1,                          Nick%Andrews%NULL%1,                          Rachel%Byford%NULL%1,                          Gavin%Dabrera%NULL%1,                          Alex%Elliot%NULL%1,                          Joanna%Ellis%NULL%0,                          Filipa%Ferreira%NULL%1,                          Jamie%Lopez Bernal%NULL%1,                          Cecilia%Okusi%NULL%1,                          Mary%Ramsay%NULL%1,                          Julian%Sherlock%NULL%1,                          Gillian%Smith%NULL%1,                          John%Williams%NULL%1,                          Gary%Howsam%NULL%1,                          Maria%Zambon%NULL%0,                          Mark%Joy%NULL%1,                          F D Richard%Hobbs%NULL%1]'
$ws.Range('E12').Value = '[Roopa%Kalyanaraman Marcello%NULL%1,                          Johanna%Dolle%NULL%2,                          Johanna%Dolle%NULL%0,                          Sheila%Grami%NULL%1,                          Richard%Adule%NULL%1,                          Zeyu%Li%NULL%1,                          Kathleen%Tatem%NULL%1,                          Chinyere%Anyaogu%NULL%1,                          Stephen%Apfelroth%NULL%1,                          Raji%Ayinla%NULL%2,                          Raji%Ayinla%NULL%0,                          Noella%Boma%NULL%1,                          Terence%Brady%NULL%1,                          Braulio F.%Cosme-Thormann%NULL%1,                          Roseann%Costarella%NULL%1,                          Kenra%Ford%NULL%1,                          Kecia%Gaither%NULL%1,                          Jessica%Jacobson%NULL%1,                          Marc%Kanter%NULL%1,                          Stuart%Kessler%NULL%1,                          Ross B.%Kristal%NULL%1,                          Joseph J.%Lieber%NULL%1,                          Vikramjit%Mukherjee%NULL%1,                          Vincent%Rizzo%NULL%1,                          Madden%Rowell%NULL%1,                          David%Stevens%NULL%1,                          Elana%Sydney%NULL%1,                          Andrew%Wallach%NULL%2,                          Andrew%Wallach%NULL%0,                          Dave A.%Chokshi%NULL%1,                          Nichola%Davis%NULL%1,                          NULL%NULL%NULL%0,                          Sze Yan%Liu%NULL%2,                          Sze Yan%Liu%NULL%0]'
$ws.Range('E13').Value = '[Sara C.%Auld%NULL%6,                          Mark%Caridi-Scheible%NULL%0,                          James M.%Blum%NULL%0,                          Chad%Robichaux%NULL%1,                          Colleen%Kraft%NULL%1,                          Jesse T.%Jacob%NULL%1,                          Craig S.%Jabaley%NULL%0,                          David%Carpenter%NULL%1,                          Roberta%Kaplow%NULL%1,                          Alfonso C.%Hernandez-Romieu%NULL%0,                          Max W.%Adelman%NULL%0,                          Greg S.%Martin%NULL%0,                          Craig M.%Coopersmith%NULL%1,                          David J.%Murphy%NULL%0,                          NULL%NULL%NULL%0,                          NULL%NULL%NULL%0]'
$ws.Range('E14').Value = '[ Diane Marie%Del Valle%null%1,                  Seunghee%Kim-Schulze%null%1,                  Hsin-Hui%Huang%null%1,                  Noam D.%Beckmann%null%1,                  Sharon%Nirenberg%null%1,                  Bo%Wang%null%1,                  Yonit%Lavin%null%1,                  Talia H.%Swartz%null%1,                  Deepu%Madduri%null%1,                  Aryeh%Stock%null%1,                  Thomas U.%Marron%null%1,                  Hui%Xie%null%1,                  Manishkumar%Patel%null%1,                  Kevin%Tuballes%null%1,                  Oliver%Oekelen%null%1,                  Adeeb%Rahman%null%1,                  Patricia%Kovatch%null%1,                  Judith A.%Aberg%null%1,                  Eric%Schadt%null%1,                  Sundar%Jagannath%null%1,                  Madhu%Mazumdar%null%1,                  Alexander W.%Charney%null%1,                  Adolfo%Firpo-Betancourt%null%1,                  Damodara Rao%Mendu%null%1,                  Jeffrey%Jhang%null%1,                  David%Reich%null%1,                  Keith%Sigel%null%1,                  Carlos%Cordon-Cardo%null%1,                  Marc%Feldmann%null%1,                  Samir%Parekh%null%1,                  Miriam%Merad%null%1,                  Sacha%Gnjatic%null%1]'
$ws.Range('E15').Value = '[Christopher A.%Martin%NULL%1,                          David R.%Jenkins%NULL%1,                          Jatinder S.%Minhas%NULL%1,                          Laura J.%Gray%NULL%1,                          Julian%Tang%NULL%1,                          Caroline%Williams%NULL%1,                          Shirley%Sze%NULL%1,                          Daniel%Pan%NULL%1,                          William%Jones%NULL%1,                          Raman%Verma%NULL%1,                          Scott%Knapp%NULL%1,                          Rupert%Major%NULL%1,                          Melanie%Davies%NULL%1,                          Nigel%Brunskill%NULL%1,                          Martin%Wiselka%NULL%1,                          Chris%Brightling%NULL%1,                          Kamlesh%Khunti%NULL%1,                          Pranab%Haldar%NULL%1,                          Manish%Pareek%mp426@le.ac.uk%1,                          NULL%NULL%NULL%0]'
$ws.Range('E16').Value = '[Thomas R%McCarty%NULL%1,                          Kelly E%Hathorn%NULL%2,                          Kelly E%Hathorn%NULL%0,                          Walker D%Redd%NULL%1,                          Nicolette J%Rodriguez%NULL%1,                          Joyce C%Zhou%NULL%1,                          Ahmad Najdat%Bazarbashi%NULL%1,                          Cheikh%Njie%NULL%1,                          Danny%Wong%NULL%1,                          Quoc-Dien%Trinh%NULL%1,                          Lin%Shen%NULL%1,                          Valerie E%Stone%NULL%2,                          Valerie E%Stone%NULL%0,                          Walter W%Chan%NULL%1]'
$ws.Range('E17').Value = '[ A.%Mendy%null%1,                           S.% Apewokin%null%1,                           A. A.% Wells%null%1,                           A. L. % Morrow%null%1,         A.%Mendy%null%1,         S.% Apewokin%null%1,         A. A.% Wells%null%1,         A. L. % Morrow%null%1]'
$ws.Range('E18').Value = '[Amy%Miles%NULL%1,                          Thomas E.%Webb%NULL%2,                          Thomas E.%Webb%NULL%0,                          Benjamin C.%Mcloughlin%NULL%2,                          Benjamin C.%Mcloughlin%NULL%0,                          Imran%Mannan%NULL%2,                          Imran%Mannan%NULL%0,                          Arshad%Rather%NULL%2,                          Arshad%Rather%NULL%0,                          Paul%Knopp%NULL%1,                          Daniel%Davis%daniel.davis@ucl.ac.uk%1]'
$ws.Range('E19').Value = '[Ana C.%Monteiro%NULL%1,                          Rajat%Suri%NULL%2,                          Rajat%Suri%NULL%0,                          Iheanacho O.%Emeruwa%NULL%2,                          Iheanacho O.%Emeruwa%NULL%0,                          Robert J.%Stretch%NULL%1,                          Roxana Y.%Cortes-Lopez%NULL%2,                          Roxana Y.%Cortes-Lopez%NULL%0,                          Alexander%Sherman%NULL%1,                          Catherine C.%Lindsay%NULL%1,                          Jennifer A.%Fulcher%NULL%1,                          David%Goodman-Meza%NULL%2,                          David%Goodman-Meza%NULL%0,                          Anil%Sapru%NULL%2,                          Anil%Sapru%NULL%0,                          Russell G.%Buhr%NULL%1,                          Steven Y.%Chang%NULL%2,                          Steven Y.%Chang%NULL%0,                          Tisha%Wang%NULL%1,                          Nida%Qadir%NULL%1,                          Giordano%Madeddu%NULL%8,                          Giordano%Madeddu%NULL%0,                          Giordano%Madeddu%NULL%0,                          Ana C.%Monteiro%NULL%2,                          Ana C.%Monteiro%NULL%0,                          Giordano%Madeddu%NULL%0,                          Giordano%Madeddu%NULL%0,                          Giordano%Madeddu%NULL%0,                          Giordano%Madeddu%NULL%0,                          Giordano%Madeddu%NULL%0]'
$ws.Range('E20').Value = '[ S.%Narain%null%1,                           D.% Stefanov%null%1,                           A. S.% Chau%null%1,                           A. G.% Weber%null%1,                           G. S.% Marder%null%1,                           B.% Kaplan%null%1,                           P.% Malhotra%null%1,                           O.% Bloom%null%1,                           A.% Liu%null%1,                           M.% Lesser%null%1,                           N. % Hajizadeh%null%1,         S.%Narain%null%1,         D.% Stefanov%null%1,         A. S.% Chau%null%1,         A. G.% Weber%null%1,         G. S.% Marder%null%1,         B.% Kaplan%null%1,         P.% Malhotra%null%1,         O.% Bloom%null%1,         A.% Liu%null%1,         M.% Lesser%null%1,         N. % Hajizadeh%null%1]'
$ws.Range('E21').Value = '[Claire L.%Niedzwiedz%NULL%1,                          Catherine A.%O’Donnell%NULL%2,                          Bhautesh Dinesh%Jani%NULL%1,                          Evangelia%Demou%NULL%1,                          Frederick K.%Ho%NULL%1,                          Carlos%Celis-Morales%NULL%1,                          Barbara I.%Nicholl%NULL%1,                          Frances S.%Mair%NULL%2,                          Paul%Welsh%NULL%2,                          Naveed%Sattar%NULL%2,                          Jill P.%Pell%NULL%2,                          S. Vittal%Katikireddi%vittal.katikireddi@glasgow.ac.uk%1]'
$ws.Range('E22').Value = '[Elizabeth J%Williamson%NULL%1,                          Alex J%Walker%NULL%1,                          Krishnan%Bhaskaran%NULL%1,                          Seb%Bacon%NULL%1,                          Chris%Bates%NULL%1,                          Caroline E%Morton%NULL%1,                          Helen J%Curtis%NULL%1,                          Amir%Mehrkar%NULL%1,                          David%Evans%NULL%1,                          Peter%Inglesby%NULL%1,                          Jonathan%Cockburn%NULL%1,                          Helen I%McDonald%NULL%1,                          Brian%MacKenna%NULL%1,                          Laurie%Tomlinson%NULL%1,                          Ian J%Douglas%NULL%1,                          Christopher T%Rentsch%NULL%1,                          Rohini%Mathur%NULL%1,                          Angel YS%Wong%NULL%1,                          Richard%Grieve%NULL%1,                          David%Harrison%NULL%1,                          Harriet%Forbes%NULL%1,                          Anna%Schultze%NULL%2,                          Richard%Croker%NULL%1,                          John%Parry%NULL%1,                          Frank%Hester%NULL%1,                          Sam%Harper%NULL%1,                          Rafael%Perera%NULL%1,                          Stephen JW%Evans%NULL%1,                          Liam%Smeeth%NULL%1,                          Ben%Goldacre%NULL%1]'
$ws.Range('E23').Value = '[Niti G.%Patel%NULL%1,                          Ajay%Bhasin%NULL%1,                          Joseph M.%Feinglass%NULL%1,                          Michael P.%Angarone%NULL%1,                          Elaine R.%Cohen%NULL%1,                          Jeffrey H.%Barsuk%NULL%1]'
$ws.Range('E24').Value = '[Kristen M. J.%Azar%xref no email%1,     Zijun%Shen%xref no email%1,     Robert J.%Romanelli%xref no email%1,     Stephen H.%Lockhart%xref no email%1,     Kelly%Smits%xref no email%1,     Sarah%Robinson%xref no email%1,     Stephanie%Brown%xref no email%1,     Alice R.%Pressman%xref no email%1]'
$ws.Range('E25').Value = '[Maulin%Patel%NULL%1,                          Andrew%Gangemi%NULL%1,                          Robert%Marron%NULL%1,                          Junad%Chowdhury%NULL%1,                          Ibraheem%Yousef%NULL%1,                          Matthew%Zheng%NULL%1,                          Nicole%Mills%NULL%1,                          Lauren%Tragesser%NULL%1,                          Julie%Giurintano%NULL%1,                          Rohit%Gupta%NULL%0,                          Matthew%Gordon%NULL%0,                          Parth%Rali%NULL%2,                          Parth%Rali%NULL%0,                          Gilbert%D''Alonso%NULL%1,                          David%Fleece%NULL%1,                          Huaqing%Zhao%NULL%0,                          Nicole%Patlakh%NULL%0,                          Gerard%Criner%NULL%1]'
$ws.Range('E26').Value = '[Russell M.%Petrak%rpetrak@midcusa.com%1,                          Nathan C.%Skorodin%nskorodin@midcusa.com%1,                          Nicholas W.%Van Hise%NULL%1,                          Robert M.%Fliegelman%NULL%1,                          Jonathan%Pinsky%NULL%1,                          Vishal%Didwania%NULL%1,                          Michael%Anderson%NULL%1,                          Melina%Diaz%NULL%1,                          Kairav%Shah%NULL%1,                          Vishnu V.%Chundi%NULL%1,                          David W.%Hines%NULL%1,                          Brian P.%Harting%NULL%1,                          Kamo%Sidwha%NULL%1,                          Brian%Yu%NULL%1,                          Paul%Brune%NULL%1,                          Anjum%Owaisi%NULL%1,                          David%Beezhold%NULL%1,                          Joseph%Kent%NULL%1,                          Dana%Vais%NULL%1,                          Alice%Han%NULL%1,                          Neethi%Gowda%NULL%1,                          Nishi%Sahgal%NULL%1,                          Jan%Silverman%NULL%1,                          Jonathan%Stake%NULL%1,                          Jenie%Nepomuceno%NULL%1,                          Renuka%Heddurshetti%NULL%1]'
$ws.Range('E27').Value = '[Pablo N%Perez-Guzman%NULL%1,                          Anna%Daunt%NULL%2,                          Anna%Daunt%NULL%0,                          Sujit%Mukherjee%NULL%1,                          Peter%Crook%NULL%1,                          Roberta%Forlano%NULL%1,                          Mara D%Kont%NULL%1,                          Alessandra%Løchen%NULL%1,                          Michaela%Vollmer%NULL%1,                          Paul%Middleton%NULL%1,                          Rebekah%Judge%NULL%1,                          Christopher%Harlow%NULL%1,                          Anet%Soubieres%NULL%1,                          Graham%Cooke%NULL%1,                          Peter J%White%NULL%1,                          Timothy B%Hallett%NULL%1,                          Paul%Aylin%NULL%1,                          Neil%Ferguson%NULL%1,                          Katharina%Hauck%NULL%1,                          Mark R%Thursz%NULL%1,                          Shevanthi%Nayagam%NULL%1]'
$ws.Range('E28').Value = '[Eboni G.%Price-Haywood%NULL%2,                          Jeffrey%Burton%NULL%2,                          Jeffrey%Burton%NULL%0,                          Daniel%Fort%NULL%1,                          Leonardo%Seoane%NULL%1]'
$ws.Range('E29').Value = '[Juliana Cepelowicz%Rajter%NULL%1,                          Michael S.%Sherman%NULL%1,                          Naaz%Fatteh%NULL%1,                          Fabio%Vogel%NULL%1,                          Jamie%Sacks%NULL%1,                          Jean-Jacques%Rajter%NULL%1]'
$ws.Range('E30').Value = '[ M.%Ramaswamy%null%1,                           P.% Mannam%null%1,                           R.% Comer%null%1,                           E.% Sinclair%null%1,                           D. B.% McQuaid%null%1,                           M. L. % Schmidt%null%1,         M.%Ramaswamy%null%1,         P.% Mannam%null%1,         R.% Comer%null%1,         E.% Sinclair%null%1,         D. B.% McQuaid%null%1,         M. L. % Schmidt%null%1]'
$ws.Range('E31').Value = '[Christopher T.%Rentsch%NULL%1,                          Farah%Kidwai-Khan%NULL%1,                          Janet P.%Tate%NULL%1,                          Lesley S.%Park%NULL%1,                          Joseph T.%King%NULL%1,                          Melissa%Skanderson%NULL%1,                          Ronald G.%Hauser%NULL%1,                          Anna%Schultze%NULL%0,                          Christopher I.%Jarvis%NULL%1,                          Mark%Holodniy%NULL%1,                          Vincent Lo%Re%NULL%1,                          Kathleen M.%Akgün%NULL%1,                          Kristina%Crothers%NULL%1,                          Tamar H.%Taddei%NULL%1,                          Matthew S.%Freiberg%NULL%1,                          Amy C.%Justice%NULL%1,                           C. T.%Rentsch%null%1,                           F.% Kidwai-Khan%null%1,                           J. P.% Tate%null%1,                           L. S.% Park%null%1,                           J. T.% King%null%1,                           M.% Skanderson%null%1,                           R. G.% Hauser%null%1,                           A.% Schultze%null%1,                           C. I.% Jarvis%null%1,                           M.% Holodniy%null%1,                           V.% Lo Re%null%1,                           K. M.% Akgun%null%1,                           K.% Crothers%null%1,                           T. H.% Taddei%null%1,                           M. S.% Freiberg%null%1,                           A. C. % Justice%null%1,                          C. T.%Rentsch%null%1,                          F.% Kidwai-Khan%null%1,                          J. P.% Tate%null%1,                          L. S.% Park%null%1,                          J. T.% King%null%1,                          M.% Skanderson%null%1,                          R. G.% Hauser%null%1,                          A.% Schultze%null%1,                          C. I.% Jarvis%null%1,                          M.% Holodniy%null%1,                          V.% Lo Re%null%1,                          K. M.% Akgun%null%1,                          K.% Crothers%null%1,                          T. H.% Taddei%null%1,                          M. S.% Freiberg%null%1,                          A. C. % Justice%null%1,         C. T.%Rentsch%null%1,         F.% Kidwai-Khan%null%1,         J. P.% Tate%null%1,         L. S.% Park%null%1,         J. T.% King%null%1,         M.% Skanderson%null%1,         R. G.% Hauser%null%1,         A.% Schultze%null%1,         C. I.% Jarvis%null%1,         M.% Holodniy%null%1,         V.% Lo Re%null%1,         K. M.% Akgun%null%1,         K.% Crothers%null%1,         T. H.% Taddei%null%1,         M. S.% Freiberg%null%1,         A. C. % Justice%null%1]'
$ws.Range('E32').Value = '[Eli S.%Rosenberg%xref no email%1,     Elizabeth M.%Dufort%xref no email%1,     Tomoko%Udo%xref no email%1,     Larissa A.%Wilberschied%xref no email%1,     Jessica%Kumar%xref no email%1,     James%Tesoriero%xref no email%1,     Patti%Weinberg%xref no email%1,     James%Kirkwood%xref no email%1,     Alison%Muse%xref no email%1,     Jack%DeHovitz%xref no email%1,     Debra S.%Blog%xref no email%1,     Brad%Hutton%xref no email%1,     David R.%Holtgrave%xref no email%1,     Howard A.%Zucker%xref no email%1]'
$ws.Range('E33').Value = '[Yelena%Rozenfeld%Yelena.Rozenfeld@providence.org%1,                          Jennifer%Beam%NULL%2,                          Jennifer%Beam%NULL%0,                          Haley%Maier%NULL%1,                          Whitney%Haggerson%NULL%1,                          Karen%Boudreau%NULL%1,                          Jamie%Carlson%NULL%1,                          Rhonda%Medows%NULL%1]'
$ws.Range('E34').Value = '[Allie%Sakowicz%NULL%1,                          Alexander E.%Ayala%NULL%1,                          Chideraa C.%Ukeje%NULL%1,                          Celeste S.%Witting%NULL%1,                          William A.%Grobman%NULL%1,                          Emily S.%Miller%NULL%1]'
$ws.Range('E35').Value = '[Simon E.%Brill%simon.brill@nhs.net%3,                          Hannah C.%Jarvis%NULL%2,                          Hannah C.%Jarvis%NULL%0,                          Ezgi%Ozcan%NULL%1,                          Thomas L. P.%Burns%NULL%1,                          Rabia A.%Warraich%NULL%1,                          Lisa J.%Amani%NULL%1,                          Amina%Jaffer%NULL%1,                          Stephanie%Paget%NULL%1,                          Anand%Sivaramakrishnan%NULL%1,                          Dean D.%Creer%NULL%1]'
$ws.Range('E36').Value = '[Elizabeth%Sapey%NULL%1,                          Suzy%Gallier%NULL%2,                          Suzy%Gallier%NULL%0,                          Chris%Mainey%NULL%1,                          Peter%Nightingale%NULL%2,                          Peter%Nightingale%NULL%0,                          David%McNulty%NULL%1,                          Hannah%Crothers%NULL%1,                          Felicity%Evison%NULL%1,                          Katharine%Reeves%NULL%1,                          Domenico%Pagano%NULL%1,                          Alastair K%Denniston%NULL%1,                          Krishnarajah%Nirantharakumar%NULL%1,                          Peter%Diggle%NULL%1,                          Simon%Ball%NULL%1,                          Lylah%Irshad%NULL%1,                          Maxim%Harris%NULL%1,                          Theodore%Nabav%NULL%1,                          A%Kolesnyk%NULL%1,                          M%Ahmed%NULL%1,                          A%Liaqat%NULL%1,                          Tanya%Pankhurst%NULL%1,                          Jamie%Coleman%NULL%1,                          Chirag%Dave%NULL%1,                          Khaled%ElFandi%NULL%1,                          Rifat%Rashid%NULL%1,                          Paul%Cockwell%NULL%1]'
$ws.Range('E37').Value = '[Emily C%Somers%NULL%1,                          Gregory A%Eschenauer%NULL%2,                          Gregory A%Eschenauer%NULL%0,                          Jonathan P%Troost%NULL%1,                          Jonathan L%Golob%NULL%1,                          Tejal N%Gandhi%NULL%1,                          Lu%Wang%NULL%1,                          Nina%Zhou%NULL%1,                          Lindsay A%Petty%NULL%1,                          Ji Hoon%Baang%NULL%1,                          Nicholas O%Dillman%NULL%1,                          David%Frame%NULL%1,                          Kevin S%Gregg%NULL%1,                          Dan R%Kaul%NULL%1,                          Jerod%Nagel%NULL%1,                          Twisha S%Patel%NULL%1,                          Shiwei%Zhou%NULL%1,                          Adam S%Lauring%NULL%1,                          David A%Hanauer%NULL%1,                          Emily%Martin%NULL%1,                          Pratima%Sharma%NULL%1,                          Christopher M%Fung%NULL%1,                          Jason M%Pogue%NULL%1]'
$ws.Range('E38').Value = '[Sara Y.%Tartof%NULL%1,                          Lei%Qian%NULL%2,                          Lei%Qian%NULL%0,                          Vennis%Hong%NULL%1,                          Rong%Wei%NULL%2,                          Rong%Wei%NULL%0,                          Ron F.%Nadjafi%NULL%2,                          Ron F.%Nadjafi%NULL%0,                          Heidi%Fischer%NULL%1,                          Zhuoxin%Li%NULL%1,                          Sally F.%Shaw%NULL%1,                          Susan L.%Caparosa%NULL%2,                          Susan L.%Caparosa%NULL%0,                          Claudia L.%Nau%NULL%1,                          Tanmai%Saxena%NULL%1,                          Gunter K.%Rieg%NULL%1,                          Bradley K.%Ackerson%NULL%1,                          Adam L.%Sharp%NULL%2,                          Adam L.%Sharp%NULL%0,                          Jacek%Skarbinski%NULL%2,                          Jacek%Skarbinski%NULL%0,                          Tej K.%Naik%NULL%1,                          Sameer B.%Murali%NULL%2,                          Sameer B.%Murali%NULL%0]'
$ws.Range('E39').Value = '[Farhaan S%Vahidy%NULL%1,                          Juan Carlos%Nicolas%NULL%2,                          Juan Carlos%Nicolas%NULL%0,                          Jennifer R%Meeks%NULL%1,                          Osman%Khan%NULL%1,                          Alan%Pan%NULL%1,                          Stephen L%Jones%NULL%1,                          Faisal%Masud%NULL%1,                          H Dirk%Sostman%NULL%1,                          Robert%Phillips%NULL%1,                          Julia D%Andrieni%NULL%1,                          Bita A%Kash%NULL%1,                          Khurram%Nasir%NULL%1]'
$ws.Range('E40').Value = '[Priya%Velu%NULL%1,                          Arryn%Craney%NULL%1,                          Phyllis%Ruggiero%NULL%1,                          John%Sipley%NULL%1,                          Lin%Cong%NULL%1,                          Erika M.%Hissong%NULL%1,                          Massimo%Loda%NULL%1,                          Lars F.%Westblade%NULL%1,                          Melissa%Cushing%NULL%2,                          Hanna%Rennert%NULL%1]'
$ws.Range('E41').Value = '[ A.-L.%WANG%null%1,                           X.% Zhong%null%1,                           Y. % Hurd%null%1,         A.-L.%WANG%null%1,         X.% Zhong%null%1,         Y. % Hurd%null%1]'
$ws.Range('E42').Value = '[Zichen%Wang%NULL%1,                          Amanda%Zheutlin%NULL%2,                          Amanda%Zheutlin%NULL%0,                          Yu-Han%Kao%NULL%1,                          Kristin%Ayers%NULL%1,                          Susan%Gross%NULL%1,                          Patricia%Kovatch%NULL%1,                          Sharon%Nirenberg%NULL%1,                          Alexander%Charney%NULL%1,                          Girish%Nadkarni%NULL%1,                          Jessica K%De Freitas%NULL%1,                          Paul%O’Reilly%NULL%1,                          Allan%Just%NULL%1,                          Carol%Horowitz%NULL%1,                          Glenn%Martin%NULL%1,                          Andrea%Branch%NULL%1,                          Benjamin S%Glicksberg%NULL%1,                          Dennis%Charney%NULL%1,                          David%Reich%NULL%1,                          William K%Oh%NULL%1,                          Eric%Schadt%NULL%1,                          Rong%Chen%NULL%0,                          Li%Li%NULL%6]'
$ws.Range('E43').Value = '[Baligh R.%Yehia%NULL%1,                          Angela%Winegar%NULL%1,                          Richard%Fogel%NULL%1,                          Mohamad%Fakih%NULL%1,                          Allison%Ottenbacher%NULL%1,                          Christine%Jesser%NULL%1,                          Angelo%Bufalino%NULL%1,                          Ren-Huai%Huang%NULL%1,                          Joseph%Cacchione%NULL%1]'
$ws.Range('E44').Value = '[Rosita%Zakeri%NULL%1,                          Rebecca%Bendayan%NULL%0,                          Mark%Ashworth%NULL%1,                          Daniel M.%Bean%NULL%1,                          Hiten%Dodhia%NULL%1,                          Stevo%Durbaba%NULL%1,                          Kevin%O''Gallagher%NULL%1,                          Claire%Palmer%NULL%1,                          Vasa%Curcin%NULL%1,                          Elizabeth%Aitken%NULL%1,                          William%Bernal%NULL%1,                          Richard D.%Barker%NULL%1,                          Sam%Norton%NULL%1,                          Martin%Gulliford%NULL%1,                          James T.H.%Teo%NULL%0,                          James%Galloway%NULL%1,                          Richard J.B.%Dobson%NULL%0,                          Ajay M.%Shah%ajay.shah@kcl.ac.uk%0]'
$ws.Range('E45').Value = '[Gabriel%Chamie%Gabriel.Chamie@ucsf.edu%1,                          Carina%Marquez%NULL%1,                          Emily%Crawford%NULL%1,                          James%Peng%NULL%1,                          Maya%Petersen%NULL%1,                          Daniel%Schwab%NULL%1,                          Joshua%Schwab%NULL%1,                          Jackie%Martinez%NULL%1,                          Diane%Jones%NULL%1,                          Douglas%Black%NULL%1,                          Monica%Gandhi%NULL%1,                          Andrew D%Kerkhoff%NULL%1,                          Vivek%Jain%NULL%1,                          Francesco%Sergi%NULL%1,                          Jon%Jacobo%NULL%1,                          Susana%Rojas%NULL%1,                          Valerie%Tulier-Laiwa%NULL%1,                          Tracy%Gallardo-Brown%NULL%1,                          Ayesha%Appa%NULL%1,                          Charles%Chiu%NULL%1,                          Mary%Rodgers%NULL%1,                          John%Hackett%NULL%1,                          NULL%NULL%NULL%0,                          Amy%Kistler%NULL%1,                          Samantha%Hao%NULL%1,                          Jack%Kamm%NULL%1,                          David%Dynerman%NULL%1,                          Joshua%Batson%NULL%1,                          Bryan%Greenhouse%NULL%1,                          Joe%DeRisi%NULL%1,                          Diane V%Havlir%NULL%1]'
$ws.Range('E46').Value = '[Andrew B.%Crouse%NULL%0,                          Tiffany%Grimes%NULL%1,                          Peng%Li%NULL%1,                          Matthew%Might%NULL%1,                          Fernando%Ovalle%NULL%1,                          Anath%Shalev%NULL%1]'
$ws.Range('E47').Value = '[Joseph E.%Ebinger%NULL%1,                          Natalie%Achamallah%NULL%1,                          Hongwei%Ji%NULL%1,                          Brian L.%Claggett%NULL%1,                          Nancy%Sun%NULL%1,                          Patrick%Botting%NULL%1,                          Trevor-Trung%Nguyen%NULL%1,                          Eric%Luong%NULL%1,                          Elizabeth H.%Kim%NULL%1,                          Eunice%Park%NULL%1,                          Yunxian%Liu%NULL%1,                          Ryan%Rosenberry%NULL%1,                          Yuri%Matusov%NULL%1,                          Steven%Zhao%NULL%1,                          Isabel%Pedraza%NULL%1,                          Tanzira%Zaman%NULL%1,                          Michael%Thompson%NULL%1,                          Koen%Raedschelders%NULL%1,                          Anders H.%Berg%NULL%2,                          Anders H.%Berg%NULL%0,                          Jonathan D.%Grein%NULL%1,                          Paul W.%Noble%NULL%1,                          Sumeet S.%Chugh%NULL%1,                          C. Noel%Bairey Merz%NULL%1,                          Eduardo%Marbán%NULL%1,                          Jennifer E.%Van Eyk%NULL%1,                          Scott D.%Solomon%NULL%1,                          Christine M.%Albert%NULL%1,                          Peter%Chen%NULL%1,                          Susan%Cheng%NULL%3,                          Yu Ru%Kou%NULL%8,                          Yu Ru%Kou%NULL%0,                          Yu Ru%Kou%NULL%0,                          Susan%Cheng%biodatacore@cshs.org%0,                          Susan%Cheng%biodatacore@cshs.org%0,                          Yu Ru%Kou%NULL%0,                          Yu Ru%Kou%NULL%0,                          Yu Ru%Kou%NULL%0,                          Yu Ru%Kou%NULL%0,                          Yu Ru%Kou%NULL%0]'
$ws.Range('E48').Value = '[Sascha%Ellington%NULL%1,                          Penelope%Strid%NULL%1,                          Van T.%Tong%NULL%1,                          Kate%Woodworth%NULL%1,                          Romeo R.%Galang%NULL%1,                          Laura D.%Zambrano%NULL%1,                          John%Nahabedian%NULL%1,                          Kayla%Anderson%NULL%1,                          Suzanne M.%Gilboa%NULL%1]'
$ws.Range('E49').Value = '[Brian T.%Garibaldi%NULL%2,                          Jacob%Fiksel%NULL%2,                          Jacob%Fiksel%NULL%0,                          John%Muschelli%NULL%2,                          John%Muschelli%NULL%0,                          Matthew L.%Robinson%NULL%2,                          Matthew L.%Robinson%NULL%0,                          Masoud%Rouhizadeh%NULL%2,                          Masoud%Rouhizadeh%NULL%0,                          Jamie%Perin%NULL%2,                          Jamie%Perin%NULL%0,                          Grant%Schumock%NULL%2,                          Grant%Schumock%NULL%0,                          Paul%Nagy%NULL%2,                          Paul%Nagy%NULL%0,                          Josh H.%Gray%NULL%2,                          Josh H.%Gray%NULL%0,                          Harsha%Malapati%NULL%2,                          Harsha%Malapati%NULL%0,                          Mariam%Ghobadi-Krueger%NULL%2,                          Mariam%Ghobadi-Krueger%NULL%0,                          Timothy M.%Niessen%NULL%1,                          Bo Soo%Kim%NULL%1,                          Peter M.%Hill%NULL%1,                          M. Shafeeq%Ahmed%NULL%1,                          Eric D.%Dobkin%NULL%1,                          Renee%Blanding%NULL%1,                          Jennifer%Abele%NULL%1,                          Bonnie%Woods%NULL%1,                          Kenneth%Harkness%NULL%1,                          David R.%Thiemann%NULL%1,                          Mary G.%Bowring%NULL%1,                          Aalok B.%Shah%NULL%2,                          Aalok B.%Shah%NULL%0,                          Mei-Cheng%Wang%NULL%1,                          Karen%Bandeen-Roche%NULL%1,                          Antony%Rosen%NULL%1,                          Scott L.%Zeger%NULL%2,                          Scott L.%Zeger%NULL%0,                          Amita%Gupta%NULL%1]'
$ws.Range('E50').Value = '[César%Caraballo%NULL%0,                          Megan%McCullough%NULL%0,                          Megan%McCullough%NULL%0,                          Michael A.%Fuery%NULL%0,                          Fouad%Chouairi%NULL%0,                          Fouad%Chouairi%NULL%0,                          Craig%Keating%NULL%0,                          Neal G.%Ravindra%NULL%0,                          P. Elliott%Miller%NULL%0,                          Maricar%Malinis%NULL%0,                          Nitu%Kashyap%NULL%0,                          Allen%Hsiao%NULL%0,                          F. Perry%Wilson%NULL%0,                          Jeptha P.%Curtis%NULL%0,                          Matthew%Grant%NULL%0,                          Eric J.%Velazquez%NULL%0,                          Nihar R.%Desai%NULL%0,                          Tariq%Ahmad%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0,                          Chiara%Lazzeri%NULL%0]'
